$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$val)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = "33.660.54"
$ws.Cells.Item(2, 5).Value = "  +6.52%  "
$ws.Cells.Item(3, 4).Value = "1.772.87"
$ws.Cells.Item(3, 5).Value = "  +3.12%  "
$ws.Cells.Item(4, 5).Value = "  +0.01%  "
Set-TextValue $ws.Cells.Item(5, 4) "223.88"
$ws.Cells.Item(5, 5).Value = "  +0.02%  "
Set-TextValue $ws.Cells.Item(6, 4) "0.557"
$ws.Cells.Item(6, 5).Value = "  +3.28%  "
$ws.Cells.Item(7, 5).Value = "  +0.10%  "
Set-TextValue $ws.Cells.Item(8, 4) "29.98"
$ws.Cells.Item(8, 5).Value = "  -0.20%  "
Set-TextValue $ws.Cells.Item(9, 4) "46.65"
$ws.Cells.Item(9, 5).Value = "  +3.86%  "
$ws.Cells.Item(10, 5).Value = "  +2.56%  "
Set-TextValue $ws.Cells.Item(11, 4) "0.0658"
$ws.Cells.Item(11, 5).Value = "  +0.69%  "
$ws.Cells.Item(12, 5).Value = "  +1.17%  "
$ws.Cells.Item(13, 4).Value = "2.027.82"
$ws.Cells.Item(13, 5).Value = "  +3.22%  "
$ws.Cells.Item(14, 4).Value = "1.774.67"
$ws.Cells.Item(14, 5).Value = "  +3.12%  "
$ws.Cells.Item(15, 5).Value = "  +0.64%  "
$ws.Cells.Item(16, 4).Value = "33.627.97"
$ws.Cells.Item(16, 5).Value = "  +6.40%  "
Set-TextValue $ws.Cells.Item(17, 4) "10.06"
$ws.Cells.Item(17, 5).Value = "  -1.07%  "
Set-TextValue $ws.Cells.Item(18, 4) "4.15"
$ws.Cells.Item(18, 5).Value = "  -1.02%  "
Set-TextValue $ws.Cells.Item(19, 4) "68.19"
Set-TextValue $ws.Cells.Item(20, 4) "249.11"
$c = $ws.Cells.Item(21, 4)
$c.Value = "0.0X0734"
$c.Characters(4,1).Text = [char]0x2083
$ws.Cells.Item(21, 5).Value = "  +1.26%  "
Set-TextValue $ws.Cells.Item(22, 4) "0.999"
$ws.Cells.Item(22, 5).Value = "  -0.08%  "
$ws.Cells.Item(24, 5).Value = "  -2.78%  "
$ws.Cells.Item(25, 5).Value = "  -2.66%  "
Set-TextValue $ws.Cells.Item(26, 4) "158.10"
$ws.Cells.Item(26, 5).Value = "  -0.67%  "
Set-TextValue $ws.Cells.Item(27, 4) "16.36"
$ws.Cells.Item(27, 5).Value = "  +1.42%  "
$ws.Cells.Item(28, 5).Value = "  +0.22%  "
$ws.Cells.Item(29, 5).Value = "  +1.24%  "
$ws.Cells.Item(30, 5).Value = "  +0.07%  "
Set-TextValue $ws.Cells.Item(31, 4) "3.78"
$ws.Cells.Item(31, 5).Value = "  -2.31%  "
Set-TextValue $ws.Cells.Item(32, 4) "0.0511"
$ws.Cells.Item(32, 5).Value = "  +1.62%  "
$ws.Cells.Item(33, 5).Value = "  +1.94%  "
Set-TextValue $ws.Cells.Item(34, 4) "3.52"
$ws.Cells.Item(34, 5).Value = "  +3.23%  "
$ws.Cells.Item(35, 5).Value = "  +3.61%  "
$ws.Cells.Item(36, 4).Value = "1.479.09"
$ws.Cells.Item(36, 5).Value = "  -3.30%  "
$ws.Cells.Item(37, 5).Value = "  +2.06%  "
$ws.Cells.Item(38, 5).Value = "  +2.06%  "
$ws.Cells.Item(39, 2).Value = "Aave"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Cells.Item(39, 4) "82.76"
$ws.Cells.Item(39, 5).Value = "  -0.50%  "
$ws.Cells.Item(40, 2).Value = "VeChain"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Cells.Item(40, 4) "0.0184"
$ws.Cells.Item(40, 5).Value = "  +1.10%  "
Set-TextValue $ws.Cells.Item(42, 4) "2.69"
$ws.Cells.Item(42, 5).Value = "  -1.89%  "
Set-TextValue $ws.Cells.Item(43, 4) "0.882"
$ws.Cells.Item(43, 5).Value = "  +3.03%  "
$ws.Cells.Item(44, 5).Value = "  +0.75%  "
Set-TextValue $ws.Cells.Item(45, 4) "0.0512"
$ws.Cells.Item(45, 5).Value = "  +1.61%  "
$ws.Cells.Item(46, 5).Value = "  +4.59%  "
$ws.Cells.Item(47, 4).Value = "1.919.33"
$ws.Cells.Item(47, 5).Value = "  +3.44%  "
$ws.Cells.Item(48, 5).Value = "  +0.07%  "
$ws.Cells.Item(49, 5).Value = "  +1.28%  "
Set-TextValue $ws.Cells.Item(50, 4) "11.64"
Set-TextValue $ws.Cells.Item(51, 4) "50.90"
$ws.Cells.Item(51, 5).Value = "  -3.41%  "
